$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with the latest scraped values.
# D-column values must stay as text (matching the sheet's existing inline-string
# cells), so we force a text number format before assigning, then clear the
# format again so no extra style gets attached to the cell.

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "67.588.26"
$dCell.ClearFormats()
$ws.Range("E2").Value = "  -0.37%  "

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "3.725.33"
$dCell.ClearFormats()
$ws.Range("E3").Value = "  -2.11%  "

$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = "0.999"
$dCell.ClearFormats()
$ws.Range("E4").Value = "  -0.06%  "

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "590.75"
$dCell.ClearFormats()
$ws.Range("E5").Value = "  -1.41%  "

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "165.29"
$dCell.ClearFormats()
$ws.Range("E6").Value = "  -2.02%  "

$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = "3.724.11"
$dCell.ClearFormats()
$ws.Range("E7").Value = "  -2.12%  "

$ws.Range("E8").Value = "  -0.04%  "

$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = "0.518"
$dCell.ClearFormats()
$ws.Range("E9").Value = "  -2.25%  "

$ws.Range("E10").Value = "  -3.82%  "

$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "6.46"
$dCell.ClearFormats()
$ws.Range("E11").Value = "  -0.76%  "

$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = "0.449"
$dCell.ClearFormats()
$ws.Range("E12").Value = "  -2.84%  "

$ws.Range("E13").Value = "  -5.47%  "

$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = "36.00"
$dCell.ClearFormats()
$ws.Range("E14").Value = "  -2.38%  "

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "4.348.64"
$dCell.ClearFormats()
$ws.Range("E15").Value = "  -2.14%  "

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "3.731.35"
$dCell.ClearFormats()
$ws.Range("E16").Value = "  -1.57%  "

$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "67.512.08"
$dCell.ClearFormats()
$ws.Range("E17").Value = "  -0.64%  "

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "18.20"
$dCell.ClearFormats()
$ws.Range("E18").Value = "  -1.04%  "

$ws.Range("E19").Value = "  -5.66%  "

$ws.Range("E20").Value = "  -0.38%  "

$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "10.67"
$dCell.ClearFormats()
$ws.Range("E21").Value = "  -1.76%  "

$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "466.40"
$dCell.ClearFormats()
$ws.Range("E22").Value = "  -0.61%  "

$ws.Range("E23").Value = "  -4.95%  "

$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "82.68"
$dCell.ClearFormats()
$ws.Range("E24").Value = "  -0.78%  "

$ws.Range("E25").Value = "  -11.11%  "

$ws.Range("E26").Value = "  -6.69%  "

$ws.Range("E27").Value = "  -2.44%  "

$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "10.02"
$dCell.ClearFormats()
$ws.Range("E28").Value = "  -2.36%  "

$ws.Range("E29").Value = "  -0.01%  "

$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "3.871.13"
$dCell.ClearFormats()
$ws.Range("E30").Value = "  -2.00%  "

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "2.76"
$dCell.ClearFormats()
$ws.Range("E31").Value = "  -5.69%  "

$ws.Range("E33").Value = "  -2.93%  "

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "29.48"
$dCell.ClearFormats()
$ws.Range("E34").Value = "  -4.25%  "

$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "9.01"
$dCell.ClearFormats()
$ws.Range("E35").Value = "  -3.37%  "

$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "3.677.02"
$dCell.ClearFormats()
$ws.Range("E36").Value = "  -2.50%  "

$ws.Range("E37").Value = "  -5.63%  "

$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "3.42"
$dCell.ClearFormats()
$ws.Range("E38").Value = "  -8.73%  "

$ws.Range("E39").Value = "  -1.68%  "

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "0.987"
$dCell.ClearFormats()
$ws.Range("E40").Value = "  -2.30%  "

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "5.73"
$dCell.ClearFormats()
$ws.Range("E41").Value = "  -4.25%  "

$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = "0.999"
$dCell.ClearFormats()
$ws.Range("E42").Value = "  -0.05%  "

$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "0.304"
$dCell.ClearFormats()
$ws.Range("E44").Value = "  -4.67%  "

$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "8.51"
$dCell.ClearFormats()
$ws.Range("E45").Value = "  -3.30%  "

$ws.Range("E46").Value = "  -3.59%  "

$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "45.19"
$dCell.ClearFormats()
$ws.Range("E47").Value = "  -2.70%  "

$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "142.94"
$dCell.ClearFormats()
$ws.Range("E48").Value = "  +0.53%  "

$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "384.25"
$dCell.ClearFormats()
$ws.Range("E49").Value = "  -6.04%  "

$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "25.34"
$dCell.ClearFormats()
$ws.Range("E50").Value = "  +0.19%  "

$ws.Range("E51").Value = "  -4.23%  "
